$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '34.722.84'
Set-TextValue 'E2' '  +1.98%  '
Set-TextValue 'D3' '1.793.97'
Set-TextValue 'E3' '  +0.02%  '
Set-TextValue 'E4' '  -0.25%  '
Set-TextValue 'D5' '223.75'
Set-TextValue 'E5' '  -1.93%  '
Set-TextValue 'D6' '0.553'
Set-TextValue 'E6' '  -0.56%  '
Set-TextValue 'D7' '0.998'
Set-TextValue 'E7' '  -0.29%  '
Set-TextValue 'D8' '32.45'
Set-TextValue 'E8' '  +4.04%  '
Set-TextValue 'D9' '0.284'
Set-TextValue 'E9' '  +0.81%  '
Set-TextValue 'D10' '0.0711'
Set-TextValue 'E10' '  +7.41%  '
Set-TextValue 'E11' '  +0.84%  '
Set-TextValue 'D12' '2.050.49'
Set-TextValue 'E12' '  -0.04%  '
Set-TextValue 'D13' '11.01'
Set-TextValue 'E13' '  -2.33%  '
Set-TextValue 'D14' '1.795.64'
Set-TextValue 'E14' '  -0.12%  '
Set-TextValue 'B15' 'Polygon'
Set-TextValue 'C15' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D15' '0.634'
Set-TextValue 'E15' '  -0.31%  '
Set-TextValue 'B16' 'WrappedBTC'
Set-TextValue 'C16' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D16' '34.706.71'
Set-TextValue 'E16' '  +1.81%  '
Set-TextValue 'D17' '4.30'
Set-TextValue 'E17' '  +1.94%  '
Set-TextValue 'D18' '69.11'
Set-TextValue 'E18' '  -0.85%  '
Set-TextValue 'D19' '254.15'
Set-TextValue 'E19' '  +0.22%  '
Set-TextValue 'D20' '0.0₃0805'
Set-TextValue 'E20' '  +8.05%  '
Set-TextValue 'D22' '10.76'
Set-TextValue 'E22' '  +2.86%  '
Set-TextValue 'D23' '4.20'
Set-TextValue 'E23' '  -2.11%  '
Set-TextValue 'E24' '  -0.21%  '
Set-TextValue 'D25' '160.33'
Set-TextValue 'E25' '  +1.49%  '
Set-TextValue 'D26' '16.37'
Set-TextValue 'E26' '  -1.59%  '
Set-TextValue 'D27' '7.12'
Set-TextValue 'E27' '  +1.21%  '
Set-TextValue 'E28' '  -0.36%  '
Set-TextValue 'E29' '  -0.26%  '
Set-TextValue 'D30' '0.0530'
Set-TextValue 'E30' '  +2.18%  '
Set-TextValue 'D31' '3.80'
Set-TextValue 'E31' '  -2.88%  '
Set-TextValue 'D33' '3.61'
Set-TextValue 'E33' '  -1.07%  '
Set-TextValue 'D34' '1.87'
Set-TextValue 'E34' '  -0.34%  '
Set-TextValue 'D35' '1.438.76'
Set-TextValue 'E35' '  -3.21%  '
Set-TextValue 'E36' '  +2.08%  '
Set-TextValue 'E37' '  -1.09%  '
Set-TextValue 'E38' '  -0.09%  '
Set-TextValue 'D39' '84.94'
Set-TextValue 'E39' '  +1.22%  '
Set-TextValue 'D40' '2.80'
Set-TextValue 'E40' '  -1.05%  '
Set-TextValue 'D41' '0.930'
Set-TextValue 'E41' '  +2.58%  '
Set-TextValue 'D42' '2.33'
Set-TextValue 'D43' '2.13'
Set-TextValue 'E43' '  +3.15%  '
Set-TextValue 'D44' '5.96'
Set-TextValue 'E44' '  +4.26%  '
Set-TextValue 'E45' '  -1.36%  '
Set-TextValue 'D46' '0.0491'
Set-TextValue 'E46' '  -4.84%  '
Set-TextValue 'D47' '1.948.02'
Set-TextValue 'E47' '  -0.10%  '
Set-TextValue 'D48' '105.59'
Set-TextValue 'E48' '  +7.58%  '
Set-TextValue 'E49' '  -0.22%  '
Set-TextValue 'E50' '  +1.00%  '
Set-TextValue 'E51' '  +7.95%  '
